$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 841234.4399999999
$ws.Range("I28").Value = 5001666
$ws.Range("K28").Value = 5001666
$ws.Range("M28").Value = -5001181
$ws.Range("H38").Value = 519.5
$ws.Range("I38").Value = 40
$ws.Range("K38").Value = 120
$ws.Range("M38").Value = 252
$ws.Range("H40").Value = 2779.625
$ws.Range("I40").Value = 2156.1667
$ws.Range("J40").Value = 4650
$ws.Range("K40").Value = 2156.1667
$ws.Range("L40").Value = 4650
$ws.Range("M40").Value = -1981.1667
$ws.Range("N40").Value = -5000
$ws.Range("H41").Value = 292.5
$ws.Range("I41").Value = 250.6
$ws.Range("J41").Value = 502
$ws.Range("K41").Value = 250.6
$ws.Range("L41").Value = 502
$ws.Range("M41").Value = 189.4
$ws.Range("N41").Value = -1382
$ws.Range("H70").Value = 1738.6
$ws.Range("I70").Value = 602
$ws.Range("J70").Value = 2022.75
$ws.Range("K70").Value = 1806
$ws.Range("L70").Value = 6068.25
$ws.Range("M70").Value = -1536
$ws.Range("N70").Value = -6608.25
$ws.Range("H73").Value = 1738.6
$ws.Range("I73").Value = 602
$ws.Range("J73").Value = 2022.75
$ws.Range("K73").Value = 1806
$ws.Range("L73").Value = 6068.25
$ws.Range("M73").Value = -870
$ws.Range("N73").Value = -7940.25
$ws.Range("H76").Value = 1117153.5
$ws.Range("I76").Value = 1670765.5
$ws.Range("J76").Value = 9929.666999999999
$ws.Range("K76").Value = 1670765.5
$ws.Range("L76").Value = 9929.666999999999
$ws.Range("M76").Value = -1670450.5
$ws.Range("N76").Value = -10559.667
$ws.Range("H79").Value = 1117153.5
$ws.Range("I79").Value = 1670765.5
$ws.Range("J79").Value = 9929.666999999999
$ws.Range("K79").Value = 1670765.5
$ws.Range("L79").Value = 9929.666999999999
$ws.Range("M79").Value = -1669673.5
$ws.Range("N79").Value = -12113.667
$ws.Range("H86").Value = 2003630.2
$ws.Range("I86").Value = 2003630.2
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2003630.2
$ws.Range("L86").Value = 0
$ws.Range("N86").Value = -2002507.2
$ws.Range("H88").Value = 1452.625
$ws.Range("I88").Value = 559.5
$ws.Range("J88").Value = 1750.3334
$ws.Range("K88").Value = 559.5
$ws.Range("L88").Value = 1750.3334
$ws.Range("M88").Value = -153.5
$ws.Range("N88").Value = -2562.3334
$ws.Range("H89").Value = 2003630.2
$ws.Range("I89").Value = 2003630.2
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 10018151
$ws.Range("L89").Value = 0
$ws.Range("N89").Value = -10012535
$ws.Range("H91").Value = 1452.625
$ws.Range("I91").Value = 559.5
$ws.Range("J91").Value = 1750.3334
$ws.Range("K91").Value = 559.5
$ws.Range("L91").Value = 1750.3334
$ws.Range("M91").Value = 844.5
$ws.Range("N91").Value = -4558.3334
$ws.Range("H98").Value = 3098.9375
$ws.Range("I98").Value = 2814.84
$ws.Range("K98").Value = 2814.84
$ws.Range("M98").Value = -1316.84
$ws.Range("H101").Value = 3440.25
$ws.Range("I101").Value = 4255.6665
$ws.Range("J101").Value = 994
$ws.Range("K101").Value = 12766.9995
$ws.Range("L101").Value = 2982
$ws.Range("M101").Value = -11144.9995
$ws.Range("N101").Value = -6226
$ws.Range("H118").Value = 629.1818
$ws.Range("I118").Value = 642.6
$ws.Range("K118").Value = 1927.8
$ws.Range("M118").Value = -270.8000000000002
$ws.Range("H122").Value = 3098.9375
$ws.Range("I122").Value = 2814.84
$ws.Range("K122").Value = 8444.52
$ws.Range("M122").Value = -5994.52
$ws.Range("H129").Value = 1082.4
$ws.Range("I129").Value = 699
$ws.Range("K129").Value = 2097
$ws.Range("M129").Value = 2903
$ws.Range("H132").Value = 4011.7292
$ws.Range("I132").Value = 3761.465
$ws.Range("K132").Value = 11284.395
$ws.Range("M132").Value = -8754.395
$ws.Range("H137").Value = 3965.5
$ws.Range("J137").Value = 6889.1665
$ws.Range("L137").Value = 20667.4995
$ws.Range("N137").Value = -25767.4995
$ws.Range("H138").Value = 5538.0415
$ws.Range("J138").Value = 8019.8
$ws.Range("L138").Value = 24059.4
$ws.Range("N138").Value = -34339.4
$ws.Range("M86").ClearContents()
$ws.Range("M89").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 7633
$ws.Range("J46").Value = 8215.571
$ws.Range("L46").Value = 8215.571
$ws.Range("N46").Value = -8853.571
$ws.Range("H102").Value = 1761.1428
$ws.Range("I102").Value = 1648.0526
$ws.Range("K102").Value = 1648.0526
$ws.Range("M102").Value = -26.05259999999998
$ws.Range("H122").Value = 4423.6294
$ws.Range("I122").Value = 3661.9302
$ws.Range("J122").Value = 7401.1816
$ws.Range("K122").Value = 10985.7906
$ws.Range("L122").Value = 22203.5448
$ws.Range("M122").Value = -8535.7906
$ws.Range("N122").Value = -27103.5448

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 11185
$ws.Range("I86").Value = 10994.8
$ws.Range("K86").Value = 10994.8
$ws.Range("M86").Value = -9871.799999999999
$ws.Range("H89").Value = 11185
$ws.Range("I89").Value = 10994.8
$ws.Range("K89").Value = 54974
$ws.Range("M89").Value = -49358
$ws.Range("H107").Value = 18758.688
$ws.Range("I107").Value = 21626.076
$ws.Range("K107").Value = 21626.076
$ws.Range("M107").Value = -19706.076
$ws.Range("H140").Value = 77834.336
$ws.Range("J140").Value = 82182.91
$ws.Range("L140").Value = 82182.91
$ws.Range("N140").Value = -92542.91

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 6439686.5
$ws.Range("I6").Value = 8193692
$ws.Range("J6").Value = 6304763
$ws.Range("K6").Value = 8193692
$ws.Range("L6").Value = 6304763
$ws.Range("M6").Value = -8193579
$ws.Range("N6").Value = -6304989
$ws.Range("H31").Value = 4179.433
$ws.Range("I31").Value = 1169.3125
$ws.Range("K31").Value = 1169.3125
$ws.Range("M31").Value = -874.3125
$ws.Range("H34").Value = 4179.433
$ws.Range("I34").Value = 1169.3125
$ws.Range("K34").Value = 1169.3125
$ws.Range("M34").Value = -967.3125
$ws.Range("H132").Value = 29458.6
$ws.Range("I132").Value = 44662
$ws.Range("K132").Value = 133986
$ws.Range("M132").Value = -131456

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 158
$ws.Range("J92").Value = 350
$ws.Range("L92").Value = 1050
$ws.Range("N92").Value = -3546
$ws.Range("H118").Value = 10847
$ws.Range("I118").Value = 8632.25
$ws.Range("K118").Value = 25896.75
$ws.Range("M118").Value = -24653.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 210.8
$ws.Range("I107").Value = 210.8
$ws.Range("K107").Value = 210.8
$ws.Range("M107").Value = 1709.2
$ws.Range("H122").Value = 4126.057
$ws.Range("I122").Value = 3940.88
$ws.Range("J122").Value = 4589
$ws.Range("K122").Value = 11822.64
$ws.Range("L122").Value = 13767
$ws.Range("M122").Value = -9372.639999999999
$ws.Range("N122").Value = -18667

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 833.3333
$ws.Range("I22").Value = 750
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 750
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -455
$ws.Range("N22").Value = -1590
$ws.Range("H27").Value = 833.3333
$ws.Range("I27").Value = 750
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 750
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -643
$ws.Range("N27").Value = -1214
$ws.Range("H46").Value = 1328.6316
$ws.Range("I46").Value = 1015.55554
$ws.Range("J46").Value = 1610.4
$ws.Range("K46").Value = 1015.55554
$ws.Range("L46").Value = 1610.4
$ws.Range("M46").Value = -827.55554
$ws.Range("N46").Value = -1986.4
$ws.Range("H55").Value = 2823.5293
$ws.Range("I55").Value = 772.8182
$ws.Range("J55").Value = 6583.1665
$ws.Range("K55").Value = 772.8182
$ws.Range("L55").Value = 6583.1665
$ws.Range("M55").Value = -599.8182
$ws.Range("N55").Value = -6929.1665
$ws.Range("H68").Value = 7496.4443
$ws.Range("I68").Value = 7819.7646
$ws.Range("K68").Value = 7819.7646
$ws.Range("M68").Value = -7070.7646
$ws.Range("H71").Value = 7496.4443
$ws.Range("I71").Value = 7819.7646
$ws.Range("K71").Value = 39098.823
$ws.Range("M71").Value = -35354.823
$ws.Range("H132").Value = 4235.909
$ws.Range("I132").Value = 3383
$ws.Range("K132").Value = 10149
$ws.Range("M132").Value = -7619
$ws.Range("H136").Value = 11627.363
$ws.Range("I136").Value = 16999.75
$ws.Range("J136").Value = 8557.429
$ws.Range("K136").Value = 50999.25
$ws.Range("L136").Value = 25672.287
$ws.Range("M136").Value = -48449.25
$ws.Range("N136").Value = -30772.287

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 14926.429
$ws.Range("I62").Value = 10289.167
$ws.Range("J62").Value = 42750
$ws.Range("K62").Value = 10289.167
$ws.Range("L62").Value = 42750
$ws.Range("M62").Value = -9665.166999999999
$ws.Range("N62").Value = -43998
$ws.Range("H65").Value = 14926.429
$ws.Range("I65").Value = 10289.167
$ws.Range("J65").Value = 42750
$ws.Range("K65").Value = 51445.835
$ws.Range("L65").Value = 213750
$ws.Range("M65").Value = -48325.835
$ws.Range("N65").Value = -219990
$ws.Range("H107").Value = 1253.75
$ws.Range("I107").Value = 1253.75
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3761.25
$ws.Range("L107").Value = 0
$ws.Range("N107").Value = -1841.25
$ws.Range("H122").Value = 4835
$ws.Range("I122").Value = 4634.2085
$ws.Range("K122").Value = 13902.6255
$ws.Range("M122").Value = -11452.6255
$ws.Range("H136").Value = 7077.185
$ws.Range("I136").Value = 7762.273
$ws.Range("J136").Value = 6606.1875
$ws.Range("K136").Value = 23286.819
$ws.Range("L136").Value = 19818.5625
$ws.Range("M136").Value = -20736.819
$ws.Range("N136").Value = -24918.5625
$ws.Range("M107").ClearContents()
